$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("Бєлобров А. О.", $true, $false, $false, $false, $false, $true, 1, $false, "Максименко А. В.", 2)
